# The workbook gains one new data row. A new record is inserted right
# before the current row 122 ("Hortaliza, Terminal La Palmera de La
# Serena - Berenjena"), which pushes the existing rows 122-183 down to
# 123-184 and extends the used range from A1:R183 to A1:R184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 122; Excel shifts rows 122-183 down
# to 123-184 (carrying their formatting/styles along automatically).
$ws.Rows.Item(122).Insert()

# Give the new row's date cell (column D) the same style as its neighbors
# (numFmt for dates) before writing its value.
$ws.Cells.Item(121,4).Copy()
$ws.Cells.Item(122,4).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 122 with the new record's data.
$ws.Cells.Item(122,1).Value2 = 8
$ws.Cells.Item(122,2).Value2 = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(122,3).Value2 = 'Coquimbo'
$ws.Cells.Item(122,4).Value2 = 44917
$ws.Cells.Item(122,5).Value2 = 4
$ws.Cells.Item(122,6).Value2 = 100112001
$ws.Cells.Item(122,7).Value2 = 'Berenjena'
$ws.Cells.Item(122,8).Value2 = 'Sin especificar'
$ws.Cells.Item(122,9).Value2 = 'Primera'
$ws.Cells.Item(122,10).Value2 = 400
$ws.Cells.Item(122,11).Value2 = 13000
$ws.Cells.Item(122,12).Value2 = 14000
$ws.Cells.Item(122,13).Value2 = 13500
$ws.Cells.Item(122,14).Value2 = '$/caja 40 unidades'
$ws.Cells.Item(122,15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(122,16).Value2 = 338
$ws.Cells.Item(122,17).Value2 = 40
$ws.Cells.Item(122,18).Value2 = 'Hortaliza'
